# unit conversion for RC Sections
# - Rename the ambiguous "SI" units label to the more descriptive "N,mm,t,s,C"
#   on the Rectangular and Tee sheets.
# - Add a new "Circular" sheet with two sample circular column sections.

$wb = $excel.ActiveWorkbook

$wsRect = $wb.Worksheets.Item("Rectangular")
$wsTee  = $wb.Worksheets.Item("Tee")

# --- Rectangular sheet: units column ---
$wsRect.Range("B2").Value = "N,mm,t,s,C"
$wsRect.Range("B3").Value = "N,mm,t,s,C"
$wsRect.Columns.Item(2).AutoFit()

# --- Tee sheet: units column ---
$wsTee.Range("B2").Value = "N,mm,t,s,C"
$wsTee.Range("B3").HorizontalAlignment = -4108
$wsTee.Columns.Item(2).AutoFit()

# --- Add the new Circular sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCirc = $wb.Worksheets.Add($null, $lastSheet)
$wsCirc.Name = "Circular"

$wsCirc.Range("A1:E1").Font.Bold = $true

$wsCirc.Range("A1").Value = "Name"
$wsCirc.Range("B1").Value = "Units"
$wsCirc.Range("C1").Value = "D"
$wsCirc.Range("D1").Value = "Cover"
$wsCirc.Range("E1").Value = "Bars1"

$wsCirc.Range("A2").Value = "Column 400dia"
$wsCirc.Range("B2").Value = "N,mm,t,s,C"
$wsCirc.Range("C2").Value = 400
$wsCirc.Range("D2").Value = 30
$wsCirc.Range("E2").Value = "10x25"

$wsCirc.Range("A3").Value = "Column 500dia, 2layers"
$wsCirc.Range("B3").Value = "N,mm,t,s,C"
$wsCirc.Range("C3").Value = 500
$wsCirc.Range("D3").Value = 30
$wsCirc.Range("E3").Value = "15x25|10x20"

$wsCirc.Range("B2:E3").HorizontalAlignment = -4108

$wsCirc.Columns.Item(1).AutoFit()
$wsCirc.Columns.Item(2).AutoFit()
$wsCirc.Columns.Item(5).AutoFit()

# --- Selections matching the author's saved view state ---
$wsRect.Range("A7").Select() | Out-Null
$wsTee.Range("B11").Select() | Out-Null
$wsCirc.Range("F16").Select() | Out-Null
